$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Fitness) holds a constant value of 7293 for rows 2..252.
# Update it to the new constant value 7310 for every data row.
$ws.Range("C2:C252").Value = 7310
